# Auto update: 2025-12-01 01:09:22
# Refresh the daily BTC/crypto decision sheet: bump the report date and
# reload the latest BTC-USD row + MACRO_SCORE figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Report date (column A, every data row shares this one string) ---
# Force the cell to stay plain text (matching how the sheet already stores
# it) instead of letting Excel auto-convert the ISO-looking string into a
# real date serial: stamp a Text format, assign the value, then drop the
# format again so the cell ends up with no explicit style, same as before.
$dateRange = $ws.Range("A2:A6")
$dateRange.NumberFormat = "@"
$dateRange.Value = "2025-12-01"
$dateRange.ClearFormats()

# --- 2. BTC-USD row (row 3) technical snapshot refresh ---
$ws.Range("D3").Value = 91492.55   # 종가 (close)
$ws.Range("E3").Value = 43.5       # RSI
$ws.Range("F3").Value = 4.75       # 5일수익률
$ws.Range("H3").Value = 40         # 3일상승확률(%)
$ws.Range("I3").Value = 46         # 5일상승확률(%)
$ws.Range("J3").Value = 50         # 10일상승확률(%)
$ws.Range("K3").Value = 56.2       # 최종점수

# --- 3. MACRO_SCORE (column N) recomputed for every row ---
$ws.Range("N2").Value = 85.87246918135976
$ws.Range("N3").Value = 85.87246918135976
$ws.Range("N4").Value = 85.87246918135976
$ws.Range("N5").Value = 85.87246918135976
$ws.Range("N6").Value = 85.87246918135976
